$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 111 and 112 have swapped match data (columns B and F:AC), while
# columns A, C, D, E stay the same (row index, Div, Div Original Name, Date).
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$row111Vals = @{}
$row112Vals = @{}
foreach ($col in $cols) {
    $row111Vals[$col] = $ws.Range("${col}111").Value2
    $row112Vals[$col] = $ws.Range("${col}112").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}111").Value2 = $row112Vals[$col]
    $ws.Range("${col}112").Value2 = $row111Vals[$col]
}
